$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")

# --- Row 2: "Published" column now records the raw published value as text
#     instead of as a boolean (handles cases where the source value isn't a
#     clean TRUE/FALSE) ---
$ws.Range("E2").Value = "True "

# --- Row 3: same treatment for the FALSE row ---
$ws.Range("E3").Value = "False "

# --- Row 4: a new customer row, added to exercise the "published value
#     varies" handling ---
$ws.Range("A4").Value = 10009439
$ws.Range("B4").Value = "Testy McTestface"
$ws.Range("C4").Value = "TN34 1RL"
$ws.Range("D4").Value = "Central Government"

# URN column keeps the same look as the URN column above it (right-aligned,
# Calibri - matches row 3's formatting).
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# CustomerName/PostCode/Sector columns pick up the Arial look used by row 2
# (rather than row 3's Calibri), then nudge WrapText off/on so the cells get
# their own explicit alignment record instead of silently reusing row 2's.
$ws.Range("B2:D2").Copy() | Out-Null
$ws.Range("B4:D4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B4:D4").WrapText = $false
